$d = $word.ActiveDocument

# ------------------------------------------------------------------
# First, gather every text-driven edit location from the Paragraphs
# collection BEFORE any table mutation (Rows.Add / Cell edits) is
# performed -- the Paragraphs index cache gets invalidated once a
# table is restructured, so do all paragraph scanning up front.
# ------------------------------------------------------------------

# 3) Two body paragraphs reference "TFS 5404." -- update the ticket
#    number to 5661 (leave the unrelated table-row mention of
#    "TFS 5404 - Allow users..." untouched).
$n = $d.Paragraphs.Count
$matches = @()
for ($i = 1; $i -le $n; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "*TFS 5404.*") {
        $matches += $i
    }
}

$replaceRanges = @()
foreach ($idx in $matches) {
    $p = $d.Paragraphs.Item($idx)
    $r = $p.Range
    $full = $r.Text
    $pos = $full.IndexOf("5404")
    if ($pos -ge 0) {
        $absStart = $r.Start + $pos
        $absEnd = $absStart + 4
        $replaceRanges += ,@($absStart, $absEnd)
    }
}

# 4) Build number reference: C36498 -> C36654
$buildRange = $null
for ($i = 1; $i -le $n; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "*C36498*") {
        $r = $p.Range
        $full = $r.Text
        $pos = $full.IndexOf("C36498")
        $absStart = $r.Start + $pos
        $absEnd = $absStart + 6
        $buildRange = @($absStart, $absEnd)
        break
    }
}

# ------------------------------------------------------------------
# Now apply the character-offset based replacements, working from the
# end of the document toward the start so earlier offsets stay valid.
# ------------------------------------------------------------------
if ($buildRange -ne $null) {
    $sub = $d.Range($buildRange[0], $buildRange[1])
    Write-Host ("Replacing build number [" + $sub.Text + "]")
    $sub.Text = "C36654"
}

for ($k = $replaceRanges.Count - 1; $k -ge 0; $k--) {
    $pair = $replaceRanges[$k]
    $sub = $d.Range($pair[0], $pair[1])
    Write-Host ("Replacing [" + $sub.Text + "] at " + $pair[0] + "-" + $pair[1])
    $sub.Text = "5661"
}

# ------------------------------------------------------------------
# 1) Title table (Table 1): replace the old TFS 5404 description with
#    the new TFS 5661 description.
# ------------------------------------------------------------------
$titleTable = $d.Tables.Item(1)
$titleCell = $titleTable.Rows.Item(1).Cells.Item(2)
$titleRange = $titleCell.Range
$titleSub = $d.Range($titleRange.Start, $titleRange.End - 1)
Write-Host ("Title before: [" + $titleSub.Text + "]")
$titleSub.Text = "TFS 5661 - Opportun!ty, Re!nforcement issue"
Write-Host ("Title after: [" + $titleSub.Text + "]")

# ------------------------------------------------------------------
# 2) Revision history table (Table 2): append a new row for the
#    2/22/2017 TFS 5661 change made by Lili Huang.
# ------------------------------------------------------------------
$historyTable = $d.Tables.Item(2)
$newRow = $historyTable.Rows.Add()
$newRow.Cells.Item(1).Range.Text = "2/22/2017"
$newRow.Cells.Item(2).Range.Text = "TFS 5661 - Opportun!ty, Re!nforcement issue"
$newRow.Cells.Item(3).Range.Text = "Lili Huang"
Write-Host ("History rows now: " + $historyTable.Rows.Count)

Write-Host "Done."
